$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6 table: switch the table's style (tableStyleId) from the custom
#    "Table_0" style to the built-in style {2C0FD269-E577-4B0A-A2B5-08F9628D1326}.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{2C0FD269-E577-4B0A-A2B5-08F9628D1326}")

# ---------------------------------------------------------------------------
# 2) Swap the deck's two themes: the design applied through the slide master
#    ("Integral") and the one used only by the notes master ("Office Theme")
#    traded places. The editable surface PowerPoint exposes for this is the
#    12-slot ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
#    reachable from the slide master's Theme - recolor it to the "Office
#    Theme" palette that used to live on the notes-master-only theme.
# ---------------------------------------------------------------------------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

$officeThemeRGB = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRGB[$i - 1]
}
